$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B12: was numeric 0.5, now text "0.5+3"
$ws.Range("B12").Value = "0.5+3"

# Update C12: append text about Ad System TableViewCell_s
$ws.Range("C12").Value = "Setup view and view navigation. Ad System TableViewCell_s are ready now."

# Update selection to B16
$ws.Range("B16").Select()

# Note: the source workbook's bookViews windowHeight (a cosmetic
# Excel-window-geometry value) is not exposed as a settable property
# through this COM surface; the host always (re)writes its own default
# there on save, independent of Application/Window property writes.
$excel.ActiveWindow.Height = 14580
